$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Row 2: Qty 4 -> 5
$ws.Range("C2").Value = 5

# Update Row 5: Qty 10 -> 1, Min price per unit 250 -> 200
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 200

# Insert two new rows of data before the existing last row (Keyboards/Logitech),
# which shifts it from row 7 down to row 9.
$ws.Rows("7:8").Insert()

# New row 7: Monitors | Asus | 1 | 231 | 500
$ws.Range("A7").Value = "Monitors"
$ws.Range("B7").Value = "Asus"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 231
$ws.Range("E7").Value = 500

# New row 8: Phones | Samsung | 5 | 500 | 900
$ws.Range("A8").Value = "Phones"
$ws.Range("B8").Value = "Samsung"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 500
$ws.Range("E8").Value = 900

# Update the selection to match the final state
$ws.Range("G12").Select()
